{"js": "// 1) The stray \"_GoBack\" bookmark (left over from the author's last edit\n//    position before this revision) is removed from its old spot (after\n//    \" the\" in the ADT/Pipe-and-Filter paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) \"Not well suited to enhancements.\" -> \"Design is not well suited for\n//    enhancements.\" in the comparison table.\nconst oldCon = context.document.body.search(\"Not well suited to enhancements.\", { matchCase: true });\noldCon.load(\"text\");\nawait context.sync();\nif (oldCon.items.length > 0) {\n  oldCon.items[0].insertText(\"Design is not well suited for enhancements.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) \"Both algorithms and data representation can be changed in individual\n//    modules without affecting others.\" -> \"Algorithms can be changed in\n//    individual modules without affecting others.\", and the \"_GoBack\"\n//    bookmark re-lands right after \"Algorithms\" (Word's marker for the last\n//    place text was edited/inserted).\nconst oldPro = context.document.body.search(\n  \"Both algorithms and data representation can be changed in individual modules without affecting others.\",\n  { matchCase: true }\n);\noldPro.load(\"text\");\nawait context.sync();\nif (oldPro.items.length > 0) {\n  oldPro.items[0].insertText(\n    \"Algorithms can be changed in individual modules without affecting others.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\nconst newWord = context.document.body.search(\"Algorithms\", { matchCase: true });\nnewWord.load(\"text\");\nawait context.sync();\nif (newWord.items.length > 0) {\n  const afterAlgorithms = newWord.items[0].getRange(Word.RangeLocation.end);\n  afterAlgorithms.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) The stray \"_GoBack\" bookmark (left over from the author's last edit\n#    position before this revision) is removed from its old spot (after\n#    \" the\" in the ADT/Pipe-and-Filter paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) \"Not well suited to enhancements.\" -> \"Design is not well suited for\n#    enhancements.\" in the comparison table.\n$rng1 = $d.Content\n$rng1.Find.Execute(\"Not well suited to enhancements.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Design is not well suited for enhancements.\", 2)\n\n# 3) \"Both algorithms and data representation can be changed in individual\n#    modules without affecting others.\" -> \"Algorithms can be changed in\n#    individual modules without affecting others.\", and the \"_GoBack\"\n#    bookmark re-lands right after \"Algorithms\" (Word's marker for the last\n#    place text was edited/inserted).\n$rng2 = $d.Content\n$rng2.Find.Execute(\"Both algorithms and data representation can be changed in individual modules without affecting others.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Algorithms can be changed in individual modules without affecting others.\", 2)\n\n$rng3 = $d.Content\n$rng3.Find.Execute(\"Algorithms\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n$rng3.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $rng3)\n"}
